$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")
$ws.Range("B2:B25").Value = 1
